$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("resultados2")
$ws.Columns("M").Insert()

$ws2 = $wb.Worksheets.Item("resultados3")
$ws2.Columns("M").Insert()
